$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E4").Value = 10.72
$ws.Range("F4").Value = 10.08

$ws.Range("D5").Value = 9.279999999999999
$ws.Range("F5").Value = 10.19

$ws.Range("D6").Value = 9.92
$ws.Range("E6").Value = 9.81
$ws.Range("H6").Value = 10.51

$ws.Range("I7").Value = 7.71

$ws.Range("F8").Value = 9.49

$ws.Range("G9").Value = 12.29
